# Auto-generated Excel COM-interop script
# Applies numeric cell updates to the Phoenix_Profits leve-profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit diff.

$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 9800
$ws.Range("I62").Value = 9000
$ws.Range("J62").Value = 13000
$ws.Range("K62").Value = 9000
$ws.Range("L62").Value = 13000
$ws.Range("M62").Value = -8376
$ws.Range("N62").Value = -14248

# Row 65
$ws.Range("H65").Value = 9800
$ws.Range("I65").Value = 9000
$ws.Range("J65").Value = 13000
$ws.Range("K65").Value = 45000
$ws.Range("L65").Value = 65000
$ws.Range("M65").Value = -41880
$ws.Range("N65").Value = -71240

# Row 70
$ws.Range("H70").Value = 2528.5
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 2528.5
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 7585.5
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -8125.5

# Row 73
$ws.Range("H73").Value = 2528.5
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 2528.5
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 7585.5
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -9457.5

# Row 86
$ws.Range("H86").Value = 5000
$ws.Range("I86").Value = 5000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 5000
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -3877

# Row 89
$ws.Range("H89").Value = 5000
$ws.Range("I89").Value = 5000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 25000
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -19384

# Row 98
$ws.Range("H98").Value = 1452.75
$ws.Range("I98").Value = 1443.2069
$ws.Range("J98").Value = 1492.2858
$ws.Range("K98").Value = 1443.2069
$ws.Range("L98").Value = 1492.2858
$ws.Range("M98").Value = 54.79310000000009
$ws.Range("N98").Value = -4488.2858

# Row 103
$ws.Range("H103").Value = 678.14813
$ws.Range("I103").Value = 347.29413
$ws.Range("J103").Value = 1240.6
$ws.Range("K103").Value = 1041.88239
$ws.Range("L103").Value = 3721.8
$ws.Range("M103").Value = -455.88239
$ws.Range("N103").Value = -4893.799999999999

# Row 122
$ws.Range("H122").Value = 1452.75
$ws.Range("I122").Value = 1443.2069
$ws.Range("J122").Value = 1492.2858
$ws.Range("K122").Value = 4329.620699999999
$ws.Range("L122").Value = 4476.857400000001
$ws.Range("M122").Value = -1879.620699999999
$ws.Range("N122").Value = -9376.857400000001

# Row 132
$ws.Range("H132").Value = 2087.723
$ws.Range("I132").Value = 2110.7888
$ws.Range("J132").Value = 1951.25
$ws.Range("K132").Value = 6332.366399999999
$ws.Range("L132").Value = 5853.75
$ws.Range("M132").Value = -3802.366399999999


# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3474.4517
$ws.Range("I32").Value = 3611.7222
$ws.Range("J32").Value = 2547.875
$ws.Range("K32").Value = 3611.7222
$ws.Range("L32").Value = 2547.875
$ws.Range("M32").Value = -3324.7222
$ws.Range("N32").Value = -3121.875

# Row 45
$ws.Range("H45").Value = 1902.2812
$ws.Range("I45").Value = 1004.087
$ws.Range("J45").Value = 4197.6665
$ws.Range("K45").Value = 1004.087
$ws.Range("L45").Value = 4197.6665
$ws.Range("M45").Value = -627.087

# Row 61
$ws.Range("H61").Value = 3107
$ws.Range("I61").Value = 2581.1035
$ws.Range("J61").Value = 4123.7334
$ws.Range("K61").Value = 2581.1035
$ws.Range("L61").Value = 4123.7334
$ws.Range("M61").Value = -2369.1035

# Row 74
$ws.Range("H74").Value = 1327.3684
$ws.Range("I74").Value = 1292.7778
$ws.Range("J74").Value = 1358.5
$ws.Range("K74").Value = 1292.7778
$ws.Range("L74").Value = 1358.5
$ws.Range("M74").Value = -418.7778000000001
$ws.Range("N74").Value = -3106.5

# Row 77
$ws.Range("H77").Value = 1327.3684
$ws.Range("I77").Value = 1292.7778
$ws.Range("J77").Value = 1358.5
$ws.Range("K77").Value = 6463.889
$ws.Range("L77").Value = 6792.5
$ws.Range("M77").Value = -2095.889
$ws.Range("N77").Value = -15528.5

# Row 132
$ws.Range("H132").Value = 2173.625
$ws.Range("I132").Value = 2072.4783
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 6217.4349
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -3687.4349

# Row 136
$ws.Range("H136").Value = 3107
$ws.Range("I136").Value = 2581.1035
$ws.Range("J136").Value = 4123.7334
$ws.Range("K136").Value = 7743.310500000001
$ws.Range("L136").Value = 12371.2002
$ws.Range("M136").Value = -5193.310500000001


# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 13889876
$ws.Range("I94").Value = 20834040
$ws.Range("J94").Value = 1549
$ws.Range("K94").Value = 20834040
$ws.Range("L94").Value = 1549
$ws.Range("M94").Value = -20833589
$ws.Range("N94").Value = -2451

# Row 122
$ws.Range("H122").Value = 69779.7
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 69779.7
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 69779.7
$ws.Range("N122").Value = -79579.7

# Row 134
$ws.Range("H134").Value = 5266.5713
$ws.Range("I134").Value = 4800.2905
$ws.Range("J134").Value = 8880.25
$ws.Range("K134").Value = 14400.8715
$ws.Range("L134").Value = 26640.75
$ws.Range("M134").Value = -11865.8715


# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4600.952
$ws.Range("I31").Value = 3387.2856
$ws.Range("J31").Value = 5207.7856
$ws.Range("K31").Value = 3387.2856
$ws.Range("L31").Value = 5207.7856
$ws.Range("M31").Value = -3092.2856
$ws.Range("N31").Value = -5797.7856

# Row 34
$ws.Range("H34").Value = 4600.952
$ws.Range("I34").Value = 3387.2856
$ws.Range("J34").Value = 5207.7856
$ws.Range("K34").Value = 3387.2856
$ws.Range("L34").Value = 5207.7856
$ws.Range("M34").Value = -3185.2856
$ws.Range("N34").Value = -5611.7856

# Row 58
$ws.Range("H58").Value = 2528.262
$ws.Range("I58").Value = 2403.3157
$ws.Range("J58").Value = 3715.25
$ws.Range("K58").Value = 2403.3157
$ws.Range("L58").Value = 3715.25
$ws.Range("M58").Value = -2200.3157

# Row 107
$ws.Range("H107").Value = 1853.6666
$ws.Range("I107").Value = 938.25
$ws.Range("J107").Value = 3684.5
$ws.Range("K107").Value = 938.25
$ws.Range("L107").Value = 3684.5
$ws.Range("M107").Value = 981.75
$ws.Range("N107").Value = -7524.5

# Row 132
$ws.Range("H132").Value = 4396
$ws.Range("I132").Value = 4520.3335
$ws.Range("J132").Value = 3650
$ws.Range("K132").Value = 13561.0005
$ws.Range("L132").Value = 10950
$ws.Range("M132").Value = -11031.0005

# Row 134
$ws.Range("H134").Value = 3152.9
$ws.Range("I134").Value = 3082
$ws.Range("J134").Value = 4500
$ws.Range("K134").Value = 9246
$ws.Range("L134").Value = 13500
$ws.Range("M134").Value = -6711

# Row 136
$ws.Range("H136").Value = 2528.262
$ws.Range("I136").Value = 2403.3157
$ws.Range("J136").Value = 3715.25
$ws.Range("K136").Value = 7209.9471
$ws.Range("L136").Value = 11145.75
$ws.Range("M136").Value = -4659.9471


# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 101
$ws.Range("H101").Value = 11282
$ws.Range("I101").Value = 11788
$ws.Range("J101").Value = 11029
$ws.Range("K101").Value = 35364
$ws.Range("L101").Value = 33087
$ws.Range("M101").Value = -32930
$ws.Range("N101").Value = -37955


# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 3621.6667
$ws.Range("I132").Value = 3376.805
$ws.Range("J132").Value = 5055.857
$ws.Range("K132").Value = 10130.415
$ws.Range("L132").Value = 15167.571
$ws.Range("M132").Value = -7600.414999999999
$ws.Range("N132").Value = -20227.571


# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 25
$ws.Range("H25").Value = 31601.4
$ws.Range("I25").Value = 38751.75
$ws.Range("J25").Value = 3000
$ws.Range("K25").Value = 38751.75
$ws.Range("L25").Value = 3000
$ws.Range("M25").Value = -38521.75
$ws.Range("N25").Value = -3460

# Row 26
$ws.Range("H26").Value = 24249.75
$ws.Range("I26").Value = 11999
$ws.Range("J26").Value = 28333.334
$ws.Range("K26").Value = 11999
$ws.Range("L26").Value = 28333.334
$ws.Range("M26").Value = -11704
$ws.Range("N26").Value = -28923.334

# Row 46
$ws.Range("H46").Value = 2238.3333
$ws.Range("I46").Value = 1734.25
$ws.Range("J46").Value = 3246.5
$ws.Range("K46").Value = 1734.25
$ws.Range("L46").Value = 3246.5
$ws.Range("M46").Value = -1546.25
$ws.Range("N46").Value = -3622.5

# Row 93
$ws.Range("H93").Value = 4274.3335
$ws.Range("I93").Value = 3655.5
$ws.Range("J93").Value = 5512
$ws.Range("K93").Value = 3655.5
$ws.Range("L93").Value = 5512
$ws.Range("M93").Value = -2407.5
$ws.Range("N93").Value = -8008

# Row 123
$ws.Range("H123").Value = 69429
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 69429
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 69429
$ws.Range("N123").Value = -79229

# Row 132
$ws.Range("H132").Value = 10148.218
$ws.Range("I132").Value = 8019.476
$ws.Range("J132").Value = 32500
$ws.Range("K132").Value = 24058.428
$ws.Range("L132").Value = 97500
$ws.Range("M132").Value = -21528.428
$ws.Range("N132").Value = -102560

# Row 136
$ws.Range("H136").Value = 9527492
$ws.Range("I136").Value = 2847.111
$ws.Range("J136").Value = 41673170
$ws.Range("K136").Value = 8541.332999999999
$ws.Range("L136").Value = 125019510
$ws.Range("M136").Value = -5991.332999999999


# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2924.6667
$ws.Range("I122").Value = 1462.5
$ws.Range("J122").Value = 4094.4
$ws.Range("K122").Value = 4387.5
$ws.Range("L122").Value = 12283.2
$ws.Range("M122").Value = -1937.5
$ws.Range("N122").Value = -17183.2

# Row 132
$ws.Range("H132").Value = 5982.3125
$ws.Range("I132").Value = 4191.0347
$ws.Range("J132").Value = 23298
$ws.Range("K132").Value = 12573.1041
$ws.Range("L132").Value = 69894
$ws.Range("M132").Value = -10043.1041


Write-Host "Applied Phoenix_Profits leve-profit updates across 8 sheets."